$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets -------------------------------
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match the page-margin convention used by the rest of the workbook
# (PageSetup margins are expressed in points - 72pt per inch).
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# Seed formatting (header style + date style) by copying from sheet1, then
# overwrite the copied values with the new sheet's own content below.
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$newSheet.Range("C1:D1").Value = $null
$ws1.Range("A1").Copy()
$newSheet.Range("C1").PasteSpecial(-4122)
$newSheet.Range("D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)

# --- 3. Header row -----------------------------------------------------------
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- 4. Data rows -------------------------------------------------------------
$data = @()
$data += ,@(45333.99999999999, 44, 18.31943775363785, 67.77236904812997)
$data += ,@(45340.99999999999, 41, 16.73238130947417, 66.48842171837286)
$data += ,@(45410.99999999999, 11, -14.83568184598799, 35.35388374269788)
$data += ,@(45417.99999999999, 8, -18.56861513664615, 33.31444674771885)
$data += ,@(45424.99999999999, 5, -20.35205739250475, 30.77240181064214)
$data += ,@(45431.99999999999, 2, -22.37133462131128, 29.3912269382346)
$data += ,@(45438.99999999999, 0, -26.93119398692653, 24.08415324128122)
$data += ,@(45445.99999999999, 0, -29.16061073848203, 20.88942400863909)
$data += ,@(45452.99999999999, 0, -30.62184448579417, 20.88749524252421)
$data += ,@(45459.99999999999, 0, -35.1000838793747, 16.28120861247835)
$data += ,@(45466.99999999999, 0, -38.04590849892423, 13.11258149236977)
$data += ,@(45473.99999999999, 0, -41.28300936680249, 9.109109770931099)
$data += ,@(45480.99999999999, 0, -43.62076999828075, 8.737718451014899)

$row = 2
foreach ($r in $data) {
    $newSheet.Cells.Item($row, 1).Value = $r[0]
    $newSheet.Cells.Item($row, 2).Value = $r[1]
    $newSheet.Cells.Item($row, 3).Value = $r[2]
    $newSheet.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Restore the original active sheet/selection state.
$ws1.Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null

Write-Host "done"
